# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" (Overview sheet) and the
# "Latest Handoff Datetime" (per-locale sheets) for the
# 8348ebcc-029b-4cc0-bd8b-e21b8e9ddaba.md file (row 7 on every table).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G7 - Latest HO Xliff Generate Date
$wsOverview.Range("G7").Value = "2016-08-20 00:47:09"

# zh-cn!H7 - Latest Handoff Datetime
$wsZhCn.Range("H7").Value = "2016-08-20 00:47:04"

# de-de!H7 - Latest Handoff Datetime
$wsDeDe.Range("H7").Value = "2016-08-20 00:47:09"
